# Update Betfair Back/Lay odds values for 2025-12-31 games (rows 2-13)
# per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("L2").Value = 1.4
$ws.Range("Q2").Value = 1.95
$ws.Range("AB2").Value = 15.5
$ws.Range("AF2").Value = 28
$ws.Range("AH2").Value = 18
$ws.Range("AM2").Value = 90

# Row 3
$ws.Range("F3").Value = 1.3
$ws.Range("G3").Value = 1.35
$ws.Range("H3").Value = 10.5
$ws.Range("I3").Value = 13.5
$ws.Range("J3").Value = 5.4
$ws.Range("K3").Value = 6.6
$ws.Range("L3").Value = 1.29
$ws.Range("M3").Value = 1.04
$ws.Range("N3").Value = 4.9
$ws.Range("P3").Value = 2.34
$ws.Range("R3").Value = 1.46
$ws.Range("S3").Value = 2.78
$ws.Range("V3").Value = 1.08
$ws.Range("W3").Value = 3.85
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 1000
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 14
$ws.Range("AD3").Value = 1000
$ws.Range("AF3").Value = 8.6
$ws.Range("AG3").Value = 11
$ws.Range("AH3").Value = 1000
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 15
$ws.Range("AL3").Value = 1000
$ws.Range("AN3").Value = 5.2

# Row 4
$ws.Range("F4").Value = 4.1
$ws.Range("G4").Value = 4.8
$ws.Range("H4").Value = 1.79
$ws.Range("I4").Value = 1.9
$ws.Range("J4").Value = 4
$ws.Range("K4").Value = 4.6
$ws.Range("N4").Value = 4.7
$ws.Range("P4").Value = 2.26
$ws.Range("Q4").Value = 1.64
$ws.Range("R4").Value = 1.5
$ws.Range("S4").Value = 2.62
$ws.Range("T4").Value = 1.63
$ws.Range("U4").Value = 2.28
$ws.Range("V4").Value = 2.1
$ws.Range("W4").Value = 1.27
$ws.Range("X4").Value = 1000
$ws.Range("Z4").Value = 1000
$ws.Range("AA4").Value = 900
$ws.Range("AB4").Value = 1000
$ws.Range("AC4").Value = 14
$ws.Range("AD4").Value = 20
$ws.Range("AK4").Value = 220
$ws.Range("AN4").Value = 140
$ws.Range("AO4").Value = 55

# Row 5
$ws.Range("F5").Value = 3.95
$ws.Range("I5").Value = 2.1
$ws.Range("J5").Value = 3.55
$ws.Range("L5").Value = 1.39
$ws.Range("V5").Value = 1.92
$ws.Range("W5").Value = 1.27
$ws.Range("Z5").Value = 13

# Row 6
$ws.Range("F6").Value = 3.8
$ws.Range("H6").Value = 1.85
$ws.Range("I6").Value = 2.02
$ws.Range("J6").Value = 3.9
$ws.Range("K6").Value = 4.6
$ws.Range("L6").Value = 1.29
$ws.Range("N6").Value = 5.1
$ws.Range("P6").Value = 2.44
$ws.Range("Q6").Value = 1.53
$ws.Range("R6").Value = 1.57
$ws.Range("S6").Value = 2.32
$ws.Range("T6").Value = 1.66
$ws.Range("U6").Value = 2.44
$ws.Range("V6").Value = 1.98
$ws.Range("W6").Value = 1.31
$ws.Range("AD6").Value = 40
$ws.Range("AG6").Value = 30
$ws.Range("AL6").Value = 420
$ws.Range("AO6").Value = 55

# Row 7
$ws.Range("G7").Value = 9.199999999999999
$ws.Range("H7").Value = 1.44
$ws.Range("J7").Value = 4.6
$ws.Range("K7").Value = 5.5
$ws.Range("L7").Value = 1.32
$ws.Range("Q7").Value = 1.7
$ws.Range("R7").Value = 1.45
$ws.Range("S7").Value = 3
$ws.Range("T7").Value = 1.9
$ws.Range("U7").Value = 1.91
$ws.Range("Y7").Value = 17
$ws.Range("Z7").Value = 16

# Row 8
$ws.Range("F8").Value = 3.25
$ws.Range("G8").Value = 3.55
$ws.Range("H8").Value = 2.06
$ws.Range("I8").Value = 2.16
$ws.Range("J8").Value = 4
$ws.Range("L8").Value = 1.26
$ws.Range("P8").Value = 2.68
$ws.Range("S8").Value = 2.26
$ws.Range("T8").Value = 1.52
$ws.Range("V8").Value = 1.86
$ws.Range("W8").Value = 1.39
$ws.Range("Z8").Value = 17.5
$ws.Range("AD8").Value = 11.5
$ws.Range("AE8").Value = 19.5

# Row 9
$ws.Range("F9").Value = 1.92
$ws.Range("G9").Value = 2.1
$ws.Range("H9").Value = 3.65
$ws.Range("M9").Value = 1.04
$ws.Range("N9").Value = 5.2
$ws.Range("P9").Value = 2.46
$ws.Range("Q9").Value = 1.58
$ws.Range("S9").Value = 2.62
$ws.Range("T9").Value = 1.55
$ws.Range("U9").Value = 2.48
$ws.Range("X9").Value = 25
$ws.Range("AD9").Value = 17
$ws.Range("AH9").Value = 16
$ws.Range("AJ9").Value = 25

# Row 10
$ws.Range("F10").Value = 5.2
$ws.Range("G10").Value = 5.9
$ws.Range("H10").Value = 1.68
$ws.Range("I10").Value = 1.74
$ws.Range("J10").Value = 4.1
$ws.Range("K10").Value = 4.5
$ws.Range("P10").Value = 2.06
$ws.Range("S10").Value = 3.05
$ws.Range("T10").Value = 1.83
$ws.Range("U10").Value = 2
$ws.Range("V10").Value = 2.34
$ws.Range("X10").Value = 30
$ws.Range("AA10").Value = 30
$ws.Range("AB10").Value = 26
$ws.Range("AC10").Value = 9.199999999999999
$ws.Range("AE10").Value = 17.5
$ws.Range("AF10").Value = 120
$ws.Range("AG10").Value = 22
$ws.Range("AH10").Value = 60
$ws.Range("AI10").Value = 38
$ws.Range("AK10").Value = 80
$ws.Range("AO10").Value = 28

# Row 11
$ws.Range("F11").Value = 14
$ws.Range("I11").Value = 1.25
$ws.Range("K11").Value = 8.6
$ws.Range("N11").Value = 8.6
$ws.Range("O11").Value = 1.1
$ws.Range("P11").Value = 3.4
$ws.Range("R11").Value = 1.98
$ws.Range("S11").Value = 1.84
$ws.Range("T11").Value = 1.81
$ws.Range("U11").Value = 2.02
$ws.Range("V11").Value = 5
$ws.Range("X11").Value = 1000
$ws.Range("Y11").Value = 16
$ws.Range("AB11").Value = 990
$ws.Range("AD11").Value = 12
$ws.Range("AG11").Value = 55
$ws.Range("AL11").Value = 150
$ws.Range("AN11").Value = 170
$ws.Range("AO11").Value = 3.2

# Row 12
$ws.Range("F12").Value = 4.5
$ws.Range("G12").Value = 5.9
$ws.Range("I12").Value = 1.78
$ws.Range("J12").Value = 4.3
$ws.Range("K12").Value = 5.3
$ws.Range("L12").Value = 1.25
$ws.Range("O12").Value = 1.21
$ws.Range("Q12").Value = 1.6
$ws.Range("S12").Value = 2.56
$ws.Range("T12").Value = 1.05
$ws.Range("U12").Value = 1.04
$ws.Range("V12").Value = 2.28
$ws.Range("W12").Value = 1.23
$ws.Range("Y12").Value = 11.5
$ws.Range("AC12").Value = 11
$ws.Range("AG12").Value = 21
$ws.Range("AK12").Value = 320
$ws.Range("AL12").Value = 160
$ws.Range("AN12").Value = 120

# Row 13
$ws.Range("F13").Value = 1.27
$ws.Range("H13").Value = 11
$ws.Range("I13").Value = 14
$ws.Range("K13").Value = 7.2
$ws.Range("N13").Value = 6.2
$ws.Range("P13").Value = 2.64
$ws.Range("R13").Value = 1.62
$ws.Range("S13").Value = 2.28
$ws.Range("T13").Value = 1.86
$ws.Range("U13").Value = 1.94
$ws.Range("AB13").Value = 12.5
$ws.Range("AG13").Value = 11
$ws.Range("AH13").Value = 990
$ws.Range("AL13").Value = 70
$ws.Range("AN13").Value = 4.5

Write-Host "Updated odds values for rows 2-13"
